# run prepare & render with final data
# Update the simulated/bootstrapped support-share values in the
# "country_comparison / main_radical_redistr_positive" sheet with the
# final values produced by the data-prep & render pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.678679804978435
$ws.Range("K2").Value = 0.67987520077397
$ws.Range("L2").Value = 0.74084626345024
$ws.Range("N2").Value = 0.640735542026883
$ws.Range("B3").Value = 0.553718303775116
$ws.Range("D3").Value = 0.62468719027894
$ws.Range("E3").Value = 0.569498424923018
$ws.Range("F3").Value = 0.759057679811368
$ws.Range("G3").Value = 0.522861429260386
$ws.Range("H3").Value = 0.68963033430472
$ws.Range("I3").Value = 0.593566149815293
$ws.Range("J3").Value = 0.610091501584142
$ws.Range("K3").Value = 0.543669469571705
$ws.Range("L3").Value = 0.489608113420232
$ws.Range("M3").Value = 0.845257807486015
$ws.Range("N3").Value = 0.476878037078152
$ws.Range("B4").Value = 0.655979457579115
$ws.Range("D4").Value = 0.717338409633504
$ws.Range("E4").Value = 0.698624076765589
$ws.Range("F4").Value = 0.817349351009749
$ws.Range("G4").Value = 0.613059109334518
$ws.Range("H4").Value = 0.753270415350046
$ws.Range("I4").Value = 0.643407610297565
$ws.Range("J4").Value = 0.67911221407656
$ws.Range("K4").Value = 0.666980642669643
$ws.Range("L4").Value = 0.586347664911907
$ws.Range("M4").Value = 0.879731656712811
$ws.Range("N4").Value = 0.597918611397897
$ws.Range("B5").Value = 0.704317439138463
$ws.Range("D5").Value = 0.754718821407338
$ws.Range("E5").Value = 0.734911648538593
$ws.Range("F5").Value = 0.824902780596881
$ws.Range("G5").Value = 0.664146429814073
$ws.Range("H5").Value = 0.732255609070291
$ws.Range("I5").Value = 0.700021689160926
$ws.Range("J5").Value = 0.620279930568098
$ws.Range("K5").Value = 0.662799981897829
$ws.Range("L5").Value = 0.73581857169661
$ws.Range("M5").Value = 0.83244969139388
$ws.Range("N5").Value = 0.655883869940124
$ws.Range("B6").Value = 0.556007222541788
$ws.Range("C6").Value = 0.612606004275779
$ws.Range("K6").Value = 0.438282499020205
$ws.Range("L6").Value = 0.595756191953926
$ws.Range("N6").Value = 0.508231644030169
$ws.Range("B7").Value = 0.50312291438834
$ws.Range("K7").Value = 0.351050882580874
$ws.Range("L7").Value = 0.604399737467109
$ws.Range("N7").Value = 0.447300643788012
$ws.Range("B8").Value = 0.680881448179833
$ws.Range("K8").Value = 0.758076861129753
$ws.Range("L8").Value = 0.688615273248795
$ws.Range("N8").Value = 0.616918649447641
$ws.Range("B9").Value = 0.609601586795904
$ws.Range("K9").Value = 0.571096670838126
$ws.Range("L9").Value = 0.778963825426238
$ws.Range("N9").Value = 0.472900191628792
$ws.Range("B10").Value = 0.675595447215337
$ws.Range("K10").Value = 0.557841849059486
$ws.Range("N10").Value = 0.666651932459956
$ws.Range("B11").Value = 0.364717906507653
$ws.Range("K11").Value = 0.22156020948145
$ws.Range("N11").Value = 0.372217577193357
$ws.Range("B12").Value = 0.347853243460036
$ws.Range("N12").Value = 0.306445646731996
$ws.Range("B13").Value = 0.410626908494325
$ws.Range("K13").Value = 0.319383802321488
$ws.Range("L13").Value = 0.389233362357354
$ws.Range("N13").Value = 0.40055514051731

Write-Output "Updated 69 cells with final run values."
